$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency price (D) and 1h volume change (E) columns with
# latest scraped values. D-column price text can look numeric (e.g. "580.50"
# or "1.00"), so force Text format before assigning to avoid Excel silently
# coercing it to a number (and dropping significant trailing digits), then
# restore the Normal style so no stray formatting is introduced.
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '66.930.02'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -0.13%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.122.10'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +0.92%  '
$ws.Range("E4").Value = '  +0.02%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '580.50'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -0.03%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '172.53'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +1.86%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("E9").Value = '  -3.42%  '
$ws.Range("E10").Value = '  -1.57%  '
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("E12").Value = '  -0.98%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '37.26'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +2.30%  '
$ws.Range("E14").Value = '  -1.05%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '3.640.05'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +0.94%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '66.916.43'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -0.01%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '7.16'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -0.49%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '3.123.82'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +1.08%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '16.31'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +0.33%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '475.65'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +2.05%  '
$ws.Range("E21").Value = '  -0.58%  '
$ws.Range("E22").Value = '  +4.94%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '83.85'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -0.31%  '
$ws.Range("E24").Value = '  +1.07%  '
$ws.Range("E25").Value = '  -2.99%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '10.34'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +2.43%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '7.92'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -1.33%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '2.36'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -1.33%  '
$ws.Range("E30").Value = '  +0.41%  '
$ws.Range("E31").Value = '  +1.09%  '
$ws.Range("E32").Value = '  -0.23%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0955'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -6.27%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("E35").Value = '  -0.66%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.975'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -3.05%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '46.91'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -1.58%  '
$ws.Range("E38").Value = '  -0.82%  '
$ws.Range("E39").Value = '  -1.30%  '
$ws.Range("E40").Value = '  -1.41%  '
$ws.Range("E41").Value = '  +1.58%  '
$ws.Range("E42").Value = '  -0.68%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '2.824.00'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +1.33%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '383.04'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +0.52%  '
$ws.Range("E45").Value = '  -1.95%  '
$ws.Range("E46").Value = '  -8.52%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '136.04'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +0.93%  '
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("E50").Value = '  -0.83%  '
